$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - headers
$ws.Range("A1").Value = "Población"
$ws.Range("B1").Value = "Municipio codigo"
$ws.Range("C1").Value = "Tipo de estudios realizados"
$ws.Range("D1").Value = "Tipo de estudios realizados, código"
$ws.Range("E1").Value = "Municipio nombre"

# Row 2
$ws.Range("A2").Value = "iaest-measure:poblacion"
$ws.Range("B2").Value = "null"
$ws.Range("C2").Value = "iaest-measure:tipo-de-estudios-realizados"
$ws.Range("D2").Value = "null"
$ws.Range("E2").Value = "sdmx-dimension:refArea"

# Row 3
$ws.Range("A3").Value = "medida"
$ws.Range("B3").Value = "null"
$ws.Range("C3").Value = "medida"
$ws.Range("D3").Value = "null"
$ws.Range("E3").Value = "dim"

# Row 4
$ws.Range("A4").Value = "xsd:double"
$ws.Range("B4").Value = "null"
$ws.Range("C4").Value = "xsd:string"
$ws.Range("D4").Value = "null"
$ws.Range("E4").Value = "URI-Municipio"
